$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "html" file-format block (eng/ara/fra) after the existing
# txt/xml/json blocks, as three new rows (11, 12, 13).

# "code" column for all three new rows
$ws.Range("A11").Value = "html"
$ws.Range("A12").Value = "html"
$ws.Range("A13").Value = "html"

# Row 11 - English
$ws.Range("C11").Value = "eng"
$ws.Range("D11").Value = $true
$ws.Range("D11").HorizontalAlignment = -4131
$ws.Range("E11").Value = "superadmin"
$ws.Range("F11").Value = "now()"

# Row 12 - Arabic
$ws.Range("C12").Value = "ara"
$ws.Range("D12").Value = $true
$ws.Range("D12").HorizontalAlignment = -4131
$ws.Range("E12").Value = "superadmin"
$ws.Range("F12").Value = "now()"

# Row 13 - French
$ws.Range("C13").Value = "fra"
$ws.Range("D13").Value = $true
$ws.Range("D13").HorizontalAlignment = -4131
$ws.Range("E13").Value = "superadmin"
$ws.Range("F13").Value = "now()"

# Translated "descr" values filled in last (ara, eng, fra order)
$ws.Range("B12").Value = "ملف html"
$ws.Range("B11").Value = "html file"
$ws.Range("B13").Value = "Fichier html"

# Leave the same trailing selection as in the source workbook: G1 through
# the bottom-right corner of the sheet (columns G..XFD, all rows).
$ws.Range("G1:XFD1048576").Select()
